$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet has columns A (index), B (WIDTH), C (HEIGHT).
# The new layout inserts two columns after the index column:
#   A (index) | B (Unnamed: 0) | C (NAME) | D (WIDTH) | E (HEIGHT)
# i.e. the old WIDTH/HEIGHT data (B,C) moves to (D,E), and new
# "Unnamed: 0" / "NAME" columns are added at B/C. A new data row
# (row 5) is also appended.

# 1) Move the existing WIDTH/HEIGHT columns (with their header style)
#    from B:C into D:E.
$ws.Range("B1:C4").Copy($ws.Range("D1:E4"))

# 2) Clear the old B:C content now that it has been duplicated to D:E.
$ws.Range("B1:C4").ClearContents()

# 3) Build the new header cells, reusing the existing bold/border header
#    style (copy style from the already-moved D1 header) before setting
#    their text.
$ws.Range("D1").Copy($ws.Range("B1"))
$ws.Range("D1").Copy($ws.Range("C1"))
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "NAME"

# 4) Fill in the "Unnamed: 0" column (mirrors the A index column) and the
#    "NAME" column values for the existing rows.
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2

$ws.Range("C2").Value = "James"
$ws.Range("C3").Value = "Timmy"
$ws.Range("C4").Value = "Sally"

# 5) Append the new 4th data row (row 5), reusing the index-column style
#    for A5.
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 3

# B5 is blank for this row (no "Unnamed: 0" value) - touch a neutral,
# already-default formatting property so the cell is materialized
# without pulling in a new style.
$ws.Range("B5").Font.Bold = $false

$ws.Range("C5").Value = "Andrew"
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 22
